# New PO forecast model
# Updates the weekly/monthly/forecast tables with the latest PO forecast run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" -- append new week row 10
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A10").Value = 45662.99999999999
$wsWeekly.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B10").Value = 1

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" -- append new month row 6
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A6").Value = 45688.99999999999
$wsMonthly.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("B6").Value = 1

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" -- refreshed forecast values for existing rows
# plus new forecasted weeks appended at the bottom
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Existing weeks: forecast quantity (column B) updated, dates unchanged
$wsForecast.Range("B2").Value = 11
$wsForecast.Range("B3").Value = 11
$wsForecast.Range("B6").Value = 10
$wsForecast.Range("B7").Value = 10
$wsForecast.Range("B8").Value = 9
$wsForecast.Range("B9").Value = 9

# Rows 10-17: both the forecast date and quantity are replaced with the
# new model's future weeks
$wsForecast.Range("A10").Value = 45662.99999999999
$wsForecast.Range("B10").Value = 7

$wsForecast.Range("A11").Value = 45669.99999999999
$wsForecast.Range("B11").Value = 7

$wsForecast.Range("A12").Value = 45676.99999999999
$wsForecast.Range("B12").Value = 6

$wsForecast.Range("A13").Value = 45683.99999999999
$wsForecast.Range("B13").Value = 6

$wsForecast.Range("A14").Value = 45690.99999999999
$wsForecast.Range("B14").Value = 6

$wsForecast.Range("A15").Value = 45697.99999999999
$wsForecast.Range("B15").Value = 6

$wsForecast.Range("A16").Value = 45704.99999999999
$wsForecast.Range("B16").Value = 6

$wsForecast.Range("A17").Value = 45711.99999999999
$wsForecast.Range("B17").Value = 6

# New row 18: one more forecast week appended
$wsForecast.Range("A18").Value = 45718.99999999999
$wsForecast.Range("B18").Value = 6

# Re-apply the date number format used by the rest of column A (style
# carries over automatically for rows 10-17 since they already had it,
# but make sure row 18 matches too).
$wsForecast.Range("A10:A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
